$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39:141 down to 40:142
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new data record
$ws.Cells.Item(39, 1).Value = 3
$ws.Cells.Item(39, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = 44459
$ws.Cells.Item(39, 5).Value = 5
$ws.Cells.Item(39, 6).Value = 100112001
$ws.Cells.Item(39, 7).Value = "Berenjena"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 125
$ws.Cells.Item(39, 11).Value = 9000
$ws.Cells.Item(39, 12).Value = 9500
$ws.Cells.Item(39, 13).Value = 9260
$ws.Cells.Item(39, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(39, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 16).Value = 154
$ws.Cells.Item(39, 17).Value = 60
$ws.Cells.Item(39, 18).Value = "Hortaliza"
